# Auto-generated Excel COM-interop script applying the Marilith_Profits price-table refresh.
# For each affected sheet/row, update the currentAveragePrice* / LevePrice* / LeveProfit* columns
# (H, I, J, K, L, M, N) to the latest market-board pull. Cells that the new pull left with no
# value are cleared so the row matches the source data exactly (no stray zero-profit columns).

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 49999
$ws.Range("J87").Value = 49999
$ws.Range("L87").Value = 49999
$ws.Range("N87").Value = -52495
$ws.Range("H90").Value = 49999
$ws.Range("J90").Value = 49999
$ws.Range("L90").Value = 149997
$ws.Range("N90").Value = -162477
$ws.Range("H116").Value = 6263.8335
$ws.Range("I116").Value = 5633.3335
$ws.Range("J116").Value = 6894.3335
$ws.Range("K116").Value = 5633.3335
$ws.Range("L116").Value = 6894.3335
$ws.Range("M116").Value = -2191.3335
$ws.Range("N116").Value = -13778.3335

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5293.8535
$ws.Range("I32").Value = 4176.2
$ws.Range("K32").Value = 4176.2
$ws.Range("M32").Value = -3889.2
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").Value = ""
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").Value = ""
$ws.Range("H88").Value = 4701.4287
$ws.Range("I88").Value = 2338.5
$ws.Range("J88").Value = 5646.6
$ws.Range("K88").Value = 2338.5
$ws.Range("L88").Value = 5646.6
$ws.Range("M88").Value = -1932.5
$ws.Range("N88").Value = -6458.6
$ws.Range("H91").Value = 4701.4287
$ws.Range("I91").Value = 2338.5
$ws.Range("J91").Value = 5646.6
$ws.Range("K91").Value = 2338.5
$ws.Range("L91").Value = 5646.6
$ws.Range("M91").Value = -934.5
$ws.Range("N91").Value = -8454.6
$ws.Range("H97").Value = 759.4
$ws.Range("I97").Value = 749.25
$ws.Range("J97").Value = 800
$ws.Range("K97").Value = 749.25
$ws.Range("L97").Value = 800
$ws.Range("M97").Value = -253.25
$ws.Range("N97").Value = -1792
$ws.Range("H110").Value = 2164
$ws.Range("I110").Value = 1455.25
$ws.Range("K110").Value = 1455.25
$ws.Range("M110").Value = 589.75
$ws.Range("H132").Value = 4438.6
$ws.Range("I132").Value = 4438.6
$ws.Range("K132").Value = 13315.8
$ws.Range("M132").Value = -10785.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 224
$ws.Range("I22").Value = 148
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 148
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = 25
$ws.Range("N22").Value = -646
$ws.Range("H86").Value = 3674.3
$ws.Range("I86").Value = 3456.3333
$ws.Range("K86").Value = 3456.3333
$ws.Range("M86").Value = -2333.3333
$ws.Range("H89").Value = 3674.3
$ws.Range("I89").Value = 3456.3333
$ws.Range("K89").Value = 17281.6665
$ws.Range("M89").Value = -11665.6665
$ws.Range("H94").Value = 2183
$ws.Range("I94").Value = 2268.375
$ws.Range("K94").Value = 2268.375
$ws.Range("M94").Value = -1817.375
$ws.Range("H134").Value = 6642.4
$ws.Range("I134").Value = 6642.4
$ws.Range("K134").Value = 19927.2
$ws.Range("M134").Value = -17392.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 461.66666
$ws.Range("I22").Value = 450
$ws.Range("J22").Value = 485
$ws.Range("K22").Value = 450
$ws.Range("L22").Value = 485
$ws.Range("M22").Value = -100
$ws.Range("N22").Value = -1185
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = ""
$ws.Range("N86").Value = ""
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = ""
$ws.Range("N89").Value = ""

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = ""
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = ""
$ws.Range("H81").Value = 6980.5
$ws.Range("J81").Value = 8674
$ws.Range("L81").Value = 26022
$ws.Range("N81").Value = -28268
$ws.Range("H84").Value = 6980.5
$ws.Range("J84").Value = 8674
$ws.Range("L84").Value = 78066
$ws.Range("N84").Value = -89298

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 429.1
$ws.Range("I2").Value = 429.1
$ws.Range("K2").Value = 429.1
$ws.Range("M2").Value = -316.1
$ws.Range("H40").Value = 9237
$ws.Range("J40").Value = 9237
$ws.Range("L40").Value = 9237
$ws.Range("N40").Value = -9539
$ws.Range("H55").Value = 22347
$ws.Range("I55").Value = 8000
$ws.Range("J55").Value = 27129.334
$ws.Range("K55").Value = 8000
$ws.Range("L55").Value = 27129.334
$ws.Range("M55").Value = -7673
$ws.Range("N55").Value = -27783.334
$ws.Range("H132").Value = 2663.3333
$ws.Range("I132").Value = 2995
$ws.Range("K132").Value = 8985
$ws.Range("M132").Value = -6455

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1353
$ws.Range("I22").Value = 255
$ws.Range("K22").Value = 255
$ws.Range("M22").Value = 40
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").Value = ""
$ws.Range("H27").Value = 1353
$ws.Range("I27").Value = 255
$ws.Range("K27").Value = 255
$ws.Range("M27").Value = -148
$ws.Range("H46").Value = 3582.2354
$ws.Range("I46").Value = 2666.6667
$ws.Range("K46").Value = 2666.6667
$ws.Range("M46").Value = -2478.6667

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 18000
$ws.Range("I45").Value = 18000
$ws.Range("K45").Value = 18000
$ws.Range("M45").Value = -17509
$ws.Range("H81").Value = 414.5
$ws.Range("I81").Value = 356.875
$ws.Range("J81").Value = 645
$ws.Range("K81").Value = 713.75
$ws.Range("L81").Value = 1290
$ws.Range("M81").Value = 347.25
$ws.Range("N81").Value = -3412
$ws.Range("H84").Value = 414.5
$ws.Range("I84").Value = 356.875
$ws.Range("J84").Value = 645
$ws.Range("K84").Value = 3568.75
$ws.Range("L84").Value = 6450
$ws.Range("M84").Value = 1735.25
$ws.Range("N84").Value = -17058
$ws.Range("H132").Value = 1965.3334
$ws.Range("I132").Value = 1358.4
$ws.Range("K132").Value = 4075.2
$ws.Range("M132").Value = -1545.2
